$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, shifting existing rows 114:148 down to 115:149
$ws.Rows.Item(114).Insert()

# Populate the new row 114 with the inserted data
$ws.Cells.Item(114, 1).Value = 10
$ws.Cells.Item(114, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(114, 3).Value = "La Araucanía"
$ws.Cells.Item(114, 4).Value = 44463
$ws.Cells.Item(114, 4).NumberFormat = $ws.Cells.Item(113, 4).NumberFormat
$ws.Cells.Item(114, 5).Value = 9
$ws.Cells.Item(114, 6).Value = 100112039
$ws.Cells.Item(114, 7).Value = "Ciboulette"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 30
$ws.Cells.Item(114, 11).Value = 6000
$ws.Cells.Item(114, 12).Value = 7000
$ws.Cells.Item(114, 13).Value = 6667
$ws.Cells.Item(114, 14).Value = "$/docena de atados"
$ws.Cells.Item(114, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(114, 16).Value = 2222
$ws.Cells.Item(114, 17).Value = 3
$ws.Cells.Item(114, 18).Value = "Hortaliza"
